# Update "想去人数" (number of people interested) figures on the
# "展览" and "全部类型" sheets, per the upstream data refresh.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F9").Value  = 2032
    $ws.Range("F10").Value = 32
    $ws.Range("F11").Value = 1179
    $ws.Range("F13").Value = 79
}
